$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a cell as plain text, avoiding Excel auto-converting
# numeric-looking strings (e.g. "2.470", "1.003") into numbers, and
# resetting the cell style afterwards so no stray number-format style
# is left behind on the cell.
function Set-CellText($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-CellText "D2" "28.577.44"
Set-CellText "E2" "  -0.72%  "

Set-CellText "D3" "1.886.84"
Set-CellText "E3" "  +0.41%  "

Set-CellText "E4" "  +0.10%  "

Set-CellText "D5" "326.91"
Set-CellText "E5" "  +0.14%  "

Set-CellText "D6" "1.007"
Set-CellText "E6" "  +0.20%  "

Set-CellText "D7" "0.4601"
Set-CellText "E7" "  -1.38%  "

Set-CellText "D8" "0.3869"
Set-CellText "E8" "  -1.96%  "

Set-CellText "D9" "46.77"
Set-CellText "E9" "  +0.41%  "

Set-CellText "D10" "0.07876"
Set-CellText "E10" "  -0.55%  "

Set-CellText "D11" "1.003"
Set-CellText "E11" "  +2.69%  "

Set-CellText "D12" "21.68"
Set-CellText "E12" "  -3.16%  "

Set-CellText "D13" "1.906.63"
Set-CellText "E13" "  +0.12%  "

Set-CellText "D14" "7.079"
Set-CellText "E14" "  +0.98%  "

Set-CellText "D15" "5.714"
Set-CellText "E15" "  -0.74%  "

Set-CellText "D16" "0.06966"
Set-CellText "E16" "  -0.21%  "

Set-CellText "D17" "87.49"
Set-CellText "E17" "  -1.26%  "

Set-CellText "D18" "1.008"
Set-CellText "E18" "  +0.14%  "

Set-CellText "D19" "0.00001004"
Set-CellText "E19" "  -0.69%  "

Set-CellText "D20" "17.18"
Set-CellText "E20" "  +1.00%  "

Set-CellText "E21" "  +0.30%  "

Set-CellText "D22" "28.606.75"
Set-CellText "E22" "  -0.79%  "

Set-CellText "D23" "5.343"
Set-CellText "E23" "  -0.37%  "

Set-CellText "D24" "11.00"
Set-CellText "E24" "  -0.92%  "

Set-CellText "D25" "2.131.88"
Set-CellText "E25" "  -0.15%  "

Set-CellText "D26" "2.058"
Set-CellText "E26" "  -2.99%  "

Set-CellText "D27" "155.00"
Set-CellText "E27" "  +0.84%  "

Set-CellText "D28" "19.36"
Set-CellText "E28" "  -0.35%  "

Set-CellText "D29" "5.848"
Set-CellText "E29" "  +1.62%  "

Set-CellText "D30" "1.952"
Set-CellText "E30" "  -2.83%  "

Set-CellText "D31" "118.53"
Set-CellText "E31" "  -1.16%  "

Set-CellText "D32" "0.09345"
Set-CellText "E32" "  -0.44%  "

Set-CellText "D33" "0.9249"
Set-CellText "E33" "  -1.87%  "

Set-CellText "D34" "5.303"
Set-CellText "E34" "  -0.36%  "

Set-CellText "D35" "1.337"
Set-CellText "E35" "  -1.16%  "

Set-CellText "D37" "0.05776"
Set-CellText "E37" "  -2.34%  "

Set-CellText "B38" "FraxShare"
Set-CellText "C38" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-CellText "D38" "7.987"
Set-CellText "E38" "  +0.59%  "

Set-CellText "B39" "TrustWalletToken"
Set-CellText "C39" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-CellText "D39" "1.149"
Set-CellText "E39" "  -0.10%  "

Set-CellText "B40" "VeChain"
Set-CellText "C40" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-CellText "D40" "0.02066"
Set-CellText "E40" "  -2.86%  "

Set-CellText "D41" "0.5675"
Set-CellText "E41" "  -0.55%  "

Set-CellText "D42" "0.1792"
Set-CellText "E42" "  -0.01%  "

Set-CellText "D43" "9.718"
Set-CellText "E43" "  -2.69%  "

Set-CellText "B44" "Decentraland"
Set-CellText "C44" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-CellText "D44" "0.5355"
Set-CellText "E44" "  +0.12%  "

Set-CellText "B45" "EnergySwap"
Set-CellText "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText "D45" "11.73"
Set-CellText "E45" "  -1.20%  "

Set-CellText "B46" "Cronos"
Set-CellText "C46" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-CellText "D46" "0.07143"
Set-CellText "E46" "  -1.47%  "

Set-CellText "D47" "2.173"
Set-CellText "E47" "  +2.27%  "

Set-CellText "B48" "NEARProtocol"
Set-CellText "C48" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-CellText "D48" "1.842"
Set-CellText "E48" "  -0.57%  "

Set-CellText "B49" "WEMIXToken"
Set-CellText "C49" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-CellText "D49" "1.117"
Set-CellText "E49" "  -3.40%  "

Set-CellText "D50" "112.51"
Set-CellText "E50" "  -1.58%  "

Set-CellText "D51" "2.470"
Set-CellText "E51" "  +4.13%  "

Write-Output "Updated cryptos list values."
